$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 25 de Marzo de 2020 a las 20:46'

$ws.Range("A6").Value = 'Estados Unidos'
$ws.Range("B6").Value = 63098
$ws.Range("C6").Value = 8242
$ws.Range("D6").Value = 393
$ws.Range("E6").Value = 61819
$ws.Range("F6").Value = 1382
$ws.Range("G6").Value = 106
$ws.Range("H6").Value = 886

$ws.Range("A7").Value = 'España'
$ws.Range("B7").Value = 47611
$ws.Range("C7").Value = 5553
$ws.Range("D7").Value = 5367
$ws.Range("E7").Value = 38799
$ws.Range("F7").Value = 2636
$ws.Range("G7").Value = 454
$ws.Range("H7").Value = 3445

$ws.Range("A8").Value = 'Alemania'
$ws.Range("B8").Value = 37323
$ws.Range("C8").Value = 4332
$ws.Range("D8").Value = 3547
$ws.Range("E8").Value = 33570
$ws.Range("F8").Value = 23
$ws.Range("G8").Value = 47
$ws.Range("H8").Value = 206

$ws.Range("A18").Value = 'Noruega'
$ws.Range("B18").Value = 3066
$ws.Range("C18").Value = 200
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 3046
$ws.Range("F18").Value = 57
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 14

$ws.Range("A33").Value = 'Pakistan'
$ws.Range("B33").Value = 1063
$ws.Range("C33").Value = 91
$ws.Range("D33").Value = 21
$ws.Range("E33").Value = 1034
$ws.Range("F33").Value = 5
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 8

$ws.Range("A34").Value = 'Polonia'
$ws.Range("B34").Value = 1031
$ws.Range("C34").Value = 130
$ws.Range("D34").Value = 2
$ws.Range("E34").Value = 1015
$ws.Range("F34").Value = 3
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 14

$ws.Range("A45").Value = 'India'
$ws.Range("B45").Value = 653
$ws.Range("C45").Value = 117
$ws.Range("D45").Value = 43
$ws.Range("E45").Value = 598
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 12

$ws.Range("A46").Value = 'Filipinas'
$ws.Range("B46").Value = 636
$ws.Range("C46").Value = 84
$ws.Range("D46").Value = 26
$ws.Range("E46").Value = 572
$ws.Range("F46").Value = 1
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 38

$ws.Range("A47").Value = 'Singapur'
$ws.Range("B47").Value = 631
$ws.Range("C47").Value = 73
$ws.Range("D47").Value = 160
$ws.Range("E47").Value = 469
$ws.Range("F47").Value = 17
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 2

$ws.Range("A52").Value = 'Egipto'
$ws.Range("B52").Value = 456
$ws.Range("C52").Value = 54
$ws.Range("D52").Value = 95
$ws.Range("E52").Value = 340
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 21

$ws.Range("A53").Value = 'Panama'
$ws.Range("B53").Value = 443
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 1
$ws.Range("E53").Value = 434
$ws.Range("F53").Value = 33
$ws.Range("G53").Value = 2
$ws.Range("H53").Value = 8

$ws.Range("A54").Value = 'Croacia'
$ws.Range("B54").Value = 442
$ws.Range("C54").Value = 60
$ws.Range("D54").Value = 22
$ws.Range("E54").Value = 419
$ws.Range("F54").Value = 6
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 1

$ws.Range("A63").Value = 'Libano'
$ws.Range("B63").Value = 333
$ws.Range("C63").Value = 15
$ws.Range("D63").Value = 20
$ws.Range("E63").Value = 307
$ws.Range("F63").Value = 4
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 6

$ws.Range("A102").Value = 'Kazajistan'
$ws.Range("B102").Value = 81
$ws.Range("C102").Value = 9
$ws.Range("D102").Value = 0
$ws.Range("E102").Value = 81
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 0

$ws.Range("A103").Value = 'Camerun'
$ws.Range("B103").Value = 75
$ws.Range("C103").Value = 9
$ws.Range("D103").Value = 2
$ws.Range("E103").Value = 72
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 1

$ws.Range("A104").Value = 'Georgia'
$ws.Range("B104").Value = 75
$ws.Range("C104").Value = 5
$ws.Range("D104").Value = 10
$ws.Range("E104").Value = 65
$ws.Range("F104").Value = 1
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 0

$ws.Range("A105").Value = 'Guadalupe'
$ws.Range("B105").Value = 73
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 0
$ws.Range("E105").Value = 72
$ws.Range("F105").Value = 4
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 1

$ws.Range("A106").Value = 'Costa de Marfil'
$ws.Range("B106").Value = 73
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 2
$ws.Range("E106").Value = 71
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 0
